$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data rows for pCa = 6, pCa = 5.8 and pCa = 5.
# These are currently rows 3, 4 and 8 in the sheet. Deleting the row for
# pCa=6 (row 3) and pCa=5.8 (row 4) first shifts everything up by two,
# so the row that used to hold pCa=5 (originally row 8) is now row 6;
# delete that one too. Whole-row deletes shift the remaining rows up,
# matching the data shuffle seen between the two versions of the sheet.
$ws.Rows("3:4").Delete()
$ws.Rows("6:6").Delete()
